$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column A (rows 3-7) first, in order, so new shared strings are appended
# in the order required to reproduce the target shared string table.
$ws.Range("A3").Value = "Description of system"
$ws.Range("A4").Value = "Use Case Diagram"
$ws.Range("A5").Value = "Class Diagrams"
$ws.Range("A6").Value = "Database design"
$ws.Range("A7").Value = "Screen Design"

# Then fill column B (rows 4-7)
$ws.Range("B4").Value = "Module to convert file (xml) to Form screen"
$ws.Range("B5").Value = "Create Android database or some other form of store"
$ws.Range("B6").Value = "Module to receive form schema"
$ws.Range("B7").Value = "Module to send completed form data"

# Finally C3, added last so its shared string index comes after all the above
$ws.Range("C3").Value = "User interface"

# Adjust column widths (engine rounds to nearest 1/6 internally, so pick values
# that land on the nearest achievable stored width to the target).
$ws.Columns.Item(1).ColumnWidth = 19.67
$ws.Columns.Item(2).ColumnWidth = 48.17

# Update the active selection to match the target end-state.
$ws.Range("C4").Select()
